$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.748.21'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.90%  '
$ws.Range('E3').Value = '  -0.85%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7772'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -6.95%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '241.44'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.001'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3154'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.43%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '25.28'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -6.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06977'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.73%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08030'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.03%  '
$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7606'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.11%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.895.44'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.46%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.251'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.56%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.82'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.05%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.813.29'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.73%  '
$ws.Range('E17').Value = '  -2.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.895'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.86%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '241.60'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.24%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007659'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.56%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.002'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.183'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +16.96%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.148.85'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.41%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1633'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.43%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.264'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '165.03'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.58'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.10%  '
$ws.Range('E29').Value = '  -2.84%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.399'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.94%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.531'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.370'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.00%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05594'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.76%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.019'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.82%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.256'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.91%  '
$ws.Range('E36').Value = '  -0.71%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.008'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.79%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.646'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.88%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01896'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.769'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.89%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4377'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.18%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '71.97'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.51%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.789'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.92%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.001'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.11%  '
$ws.Range('E45').Value = '  -1.11%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '101.90'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.51%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.018.22'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.71%  '
$ws.Range('E48').Value = '  -2.76%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.820'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.74%  '
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.070.00'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.58%  '
$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.359'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.62%  '
